# Cantines import template: rename the "siret_livreur_repas" column header
# to "siret_cuisine_centrale" (#6261).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "siret_cuisine_centrale"

# Match the author's final cursor position recorded in the saved view.
$ws.Range("C2").Select() | Out-Null
